# Genetics.xlsx - "update core follow pom" edit
#
# 1) Sheet1!A2 holds a generated "CA-xxxxxxxx" case id — it was regenerated
#    to a new random-looking id (CA-JRS0KDBT -> CA-O3SDOHZ9).
# 2) The "Pass" value in the D column (rows 7-15, 18-21, 24-27 — the
#    "Result" column of the three checklists) was cleared out, leaving the
#    cells blank but keeping their existing formatting/style.
# 3) The sheet's on-screen view was scrolled/selected differently: the
#    frozen "topLeftCell" scroll position is gone and the active
#    cell/selection moved from G20 to H12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")
$ws.Activate() | Out-Null

# Regenerate the case id in A2.
$ws.Range("A2").Value = "CA-O3SDOHZ9"

# Clear the "Pass" markers in column D for the three checklist blocks,
# keeping cell formatting intact.
$ws.Range("D7:D15").ClearContents() | Out-Null
$ws.Range("D18:D21").ClearContents() | Out-Null
$ws.Range("D24:D27").ClearContents() | Out-Null

# Update the view: scroll back to the top and move the selection to H12.
$ws.Range("A1").Select() | Out-Null
$ws.Range("H12").Select() | Out-Null
